# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" summary block (rows 3-4) ---
$ws.Range("C3").Value = 72
$ws.Range("D3").Value = 95
$ws.Range("C4").Value = 72

# --- "Good Drivers" block (rows 12-17) ---
# The rows shifted up by one driver-version slot; refresh adapter name,
# sample count, roaming % and vintage date for each row.

# Row 12: now 21.60.2.1 (was 23.100.0.4), 99.9% bucket -> 100% bucket, no vintage
$ws.Range("A12").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B12").Value = 56018
$ws.Range("D12").Value = 100
$ws.Range("E12").Value = ""

# Row 13: now 22.50.1.1 (was 22.80.0.9), 99.9% bucket -> 100% bucket, no vintage
$ws.Range("A13").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B13").Value = 34244
$ws.Range("D13").Value = 100
$ws.Range("E13").Value = ""

# Row 14: now 23.100.0.4 (was 22.50.1.1), 100% bucket -> 99.9% bucket, vintage 2024-11-10
$ws.Range("A14").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B14").Value = 442178
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").Value = "'2024-11-10"

# Row 15: now 22.80.0.9 (was 21.110.3.2), 100% bucket -> 99.9% bucket, vintage 2021-08-18
$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B15").Value = 77849
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").Value = "'2021-08-18"

# Row 16: now 21.110.3.2 (was 21.70.0.6), stays 100% bucket, vintage 2020-08-05
$ws.Range("A16").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B16").Value = 59673
$ws.Range("E16").Value = "'2020-08-05"

# Row 17: now 21.70.0.6 (was 21.60.2.1), stays 100% bucket, vintage unchanged 2019-12-14
$ws.Range("A17").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B17").Value = 113652
